$wb = $excel.ActiveWorkbook

# Hyperlink target URLs (external), matching the existing Source/Handoff file links.
$mdUrl_dedb    = "https://github.com/OpenLocalizationTest/oltest/blob/593c9aff7438ff881575552994205808e22e177c/e2e/dedb3467-9493-4ab0-8961-60a966699708.md"
$xlfUrl_zhcn   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/122a7b717b31a1ce68e783a9b2098a975412d9be/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/dedb3467-9493-4ab0-8961-60a966699708.ea5a5ce460d2cdef84e4e3672d5db72352d42d69.zh-cn.xlf"
$xlfUrl_dede   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6b6f6ade3e714023dc44236574fac905ce25139/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/dedb3467-9493-4ab0-8961-60a966699708.ea5a5ce460d2cdef84e4e3672d5db72352d42d69.de-de.xlf"

$mdFileName_dedb = "dedb3467-9493-4ab0-8961-60a966699708.md"
$xlfFileName_zhcn = "dedb3467-9493-4ab0-8961-60a966699708.ea5a5ce460d2cdef84e4e3672d5db72352d42d69.zh-cn.xlf"
$xlfFileName_dede = "dedb3467-9493-4ab0-8961-60a966699708.ea5a5ce460d2cdef84e4e3672d5db72352d42d69.de-de.xlf"

$statusText = "Handed back: in sync with en-US"

# Hyperlink font color used by the workbook's custom "HyperLink" style (RGB 100,149,237 = #6495ED)
$hlColor = 15570276

function Set-HandbackLink($ws, $cellRef, $url, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
    $f = $ws.Range($cellRef).Font
    $f.Underline = 2
    $f.Color = $hlColor
}

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------- zh-cn sheet ----------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Range("B2").Value = $statusText
Set-HandbackLink $wsZhCn "E2" $mdUrl_dedb $mdFileName_dedb
Set-HandbackLink $wsZhCn "F2" $xlfUrl_zhcn $xlfFileName_zhcn
$wsZhCn.Range("G2").Value = "2016-01-25 09:15:22"

$wsZhCn.Range("B3").Value = $statusText
Set-HandbackLink $wsZhCn "E3" $mdUrl_dedb $mdFileName_dedb
Set-HandbackLink $wsZhCn "F3" $xlfUrl_zhcn $xlfFileName_zhcn
$wsZhCn.Range("G3").Value = "2016-01-25 09:15:22"

# ---------------- de-de sheet ----------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Range("B2").Value = $statusText
Set-HandbackLink $wsDeDe "E2" $mdUrl_dedb $mdFileName_dedb
Set-HandbackLink $wsDeDe "F2" $xlfUrl_dede $xlfFileName_dede
$wsDeDe.Range("G2").Value = "2016-01-25 09:15:40"

$wsDeDe.Range("B3").Value = $statusText
Set-HandbackLink $wsDeDe "E3" $mdUrl_dedb $mdFileName_dedb
Set-HandbackLink $wsDeDe "F3" $xlfUrl_dede $xlfFileName_dede
$wsDeDe.Range("G3").Value = "2016-01-25 09:15:40"

Write-Output "done"
